$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 209.95
$ws.Range("I33").Value = 161.70589
$ws.Range("K33").Value = 161.70589
$ws.Range("M33").Value = 67.29410999999999
$ws.Range("H118").Value = 415.15384
$ws.Range("I118").Value = 270.81818
$ws.Range("J118").Value = 1209
$ws.Range("K118").Value = 812.45454
$ws.Range("L118").Value = 3627
$ws.Range("M118").Value = 844.54546
$ws.Range("N118").Value = -6941
$ws.Range("H121").Value = 1583.3125
$ws.Range("J121").Value = 1614.9788
$ws.Range("L121").Value = 4844.936400000001
$ws.Range("N121").Value = -8338.936400000001
$ws.Range("H138").Value = 3972.5876
$ws.Range("I138").Value = 2957.1177
$ws.Range("J138").Value = 4188.375
$ws.Range("K138").Value = 8871.3531
$ws.Range("L138").Value = 12565.125
$ws.Range("M138").Value = -3731.3531
$ws.Range("N138").Value = -22845.125

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22746.145
$ws.Range("I32").Value = 17150.76
$ws.Range("J32").Value = 78700
$ws.Range("K32").Value = 17150.76
$ws.Range("L32").Value = 78700
$ws.Range("M32").Value = -16863.76
$ws.Range("N32").Value = -79274
$ws.Range("H37").Value = 10587.714
$ws.Range("J37").Value = 10587.714
$ws.Range("L37").Value = 10587.714
$ws.Range("N37").Value = -11133.714
$ws.Range("H61").Value = 5850774
$ws.Range("H74").Value = 26324946
$ws.Range("I74").Value = 41667850
$ws.Range("J74").Value = 22828.572
$ws.Range("K74").Value = 41667850
$ws.Range("L74").Value = 22828.572
$ws.Range("M74").Value = -41666976
$ws.Range("N74").Value = -24576.572
$ws.Range("H77").Value = 26324946
$ws.Range("I77").Value = 41667850
$ws.Range("J77").Value = 22828.572
$ws.Range("K77").Value = 208339250
$ws.Range("L77").Value = 114142.86
$ws.Range("M77").Value = -208334882
$ws.Range("N77").Value = -122878.86
$ws.Range("H102").Value = 2410
$ws.Range("I102").Value = 2378.889
$ws.Range("J102").Value = 2550
$ws.Range("K102").Value = 2378.889
$ws.Range("L102").Value = 2550
$ws.Range("M102").Value = -756.8890000000001
$ws.Range("N102").Value = -5794
$ws.Range("H114").Value = 25000
$ws.Range("J114").Value = 25000
$ws.Range("L114").Value = 25000
$ws.Range("N114").Value = -33678
$ws.Range("H127").Value = 51660
$ws.Range("J127").Value = 51660
$ws.Range("L127").Value = 51660
$ws.Range("N127").Value = -61580
$ws.Range("H132").Value = 2158763
$ws.Range("I132").Value = 3138992
$ws.Range("K132").Value = 9416976
$ws.Range("M132").Value = -9414446
$ws.Range("H136").Value = 5850774

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8335436
$ws.Range("I134").Value = 9093022
$ws.Range("J134").Value = 1995
$ws.Range("K134").Value = 27279066
$ws.Range("L134").Value = 5985
$ws.Range("M134").Value = -27276531
$ws.Range("N134").Value = -11055

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23256.291
$ws.Range("I31").Value = 102710
$ws.Range("K31").Value = 102710
$ws.Range("M31").Value = -102415
$ws.Range("H34").Value = 23256.291
$ws.Range("I34").Value = 102710
$ws.Range("K34").Value = 102710
$ws.Range("M34").Value = -102508
$ws.Range("H50").Value = 10835.2
$ws.Range("J50").Value = 10835.2
$ws.Range("L50").Value = 10835.2
$ws.Range("N50").Value = -12085.2
$ws.Range("H51").Value = 10097.5
$ws.Range("J51").Value = 10421.875
$ws.Range("L51").Value = 10421.875
$ws.Range("N51").Value = -11893.875
$ws.Range("H58").Value = 3400.2205
$ws.Range("I58").Value = 1314.6842
$ws.Range("J58").Value = 4390.85
$ws.Range("K58").Value = 1314.6842
$ws.Range("L58").Value = 4390.85
$ws.Range("M58").Value = -1111.6842
$ws.Range("N58").Value = -4796.85
$ws.Range("H61").Value = 10097.5
$ws.Range("J61").Value = 10421.875
$ws.Range("L61").Value = 10421.875
$ws.Range("N61").Value = -11117.875
$ws.Range("H62").Value = 2868.077
$ws.Range("I62").Value = 2940.7144
$ws.Range("J62").Value = 2783.3333
$ws.Range("K62").Value = 2940.7144
$ws.Range("L62").Value = 2783.3333
$ws.Range("M62").Value = -2316.7144
$ws.Range("N62").Value = -4031.3333
$ws.Range("H65").Value = 2868.077
$ws.Range("I65").Value = 2940.7144
$ws.Range("J65").Value = 2783.3333
$ws.Range("K65").Value = 14703.572
$ws.Range("L65").Value = 13916.6665
$ws.Range("M65").Value = -11583.572
$ws.Range("N65").Value = -20156.6665
$ws.Range("H68").Value = 18530.25
$ws.Range("J68").Value = 19139.143
$ws.Range("L68").Value = 19139.143
$ws.Range("N68").Value = -20637.143
$ws.Range("H71").Value = 18530.25
$ws.Range("J71").Value = 19139.143
$ws.Range("L71").Value = 57417.429
$ws.Range("N71").Value = -64905.429
$ws.Range("H74").Value = 14381.286
$ws.Range("J74").Value = 18059.8
$ws.Range("L74").Value = 18059.8
$ws.Range("N74").Value = -19807.8
$ws.Range("H77").Value = 14381.286
$ws.Range("J77").Value = 18059.8
$ws.Range("L77").Value = 54179.39999999999
$ws.Range("N77").Value = -62915.39999999999
$ws.Range("H93").Value = 15702
$ws.Range("I93").Value = 13319
$ws.Range("K93").Value = 13319
$ws.Range("M93").Value = -11447
$ws.Range("H112").Value = 11567.333
$ws.Range("J112").Value = 14702
$ws.Range("L112").Value = 14702
$ws.Range("N112").Value = -17656
$ws.Range("H132").Value = 2572
$ws.Range("I132").Value = 2178.6667
$ws.Range("J132").Value = 2808
$ws.Range("K132").Value = 6536.000100000001
$ws.Range("L132").Value = 8424
$ws.Range("M132").Value = -4006.000100000001
$ws.Range("N132").Value = -13484
$ws.Range("H134").Value = 2373.6736
$ws.Range("I134").Value = 2110.0967
$ws.Range("J134").Value = 2827.611
$ws.Range("K134").Value = 6330.2901
$ws.Range("L134").Value = 8482.832999999999
$ws.Range("M134").Value = -3795.2901
$ws.Range("N134").Value = -13552.833
$ws.Range("H136").Value = 3400.2205
$ws.Range("I136").Value = 1314.6842
$ws.Range("J136").Value = 4390.85
$ws.Range("K136").Value = 3944.0526
$ws.Range("L136").Value = 13172.55
$ws.Range("M136").Value = -1394.0526
$ws.Range("N136").Value = -18272.55

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1396.3334
$ws.Range("I110").Value = 1396.3334
$ws.Range("K110").Value = 4189.0002
$ws.Range("M110").Value = -99.0002000000004
$ws.Range("H131").Value = 30705342
$ws.Range("J131").Value = 14708126
$ws.Range("L131").Value = 44124378
$ws.Range("N131").Value = -44134458
$ws.Range("H140").Value = 1932.7587
$ws.Range("I140").Value = 1452.5
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 4357.5
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 822.5
$ws.Range("N140").Value = -19360

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 43333
$ws.Range("J101").Value = 43333
$ws.Range("L101").Value = 43333
$ws.Range("N101").Value = -49823
$ws.Range("H102").Value = 1147.619
$ws.Range("I102").Value = 1087.5
$ws.Range("J102").Value = 1340
$ws.Range("K102").Value = 1087.5
$ws.Range("L102").Value = 1340
$ws.Range("M102").Value = 534.5
$ws.Range("N102").Value = -4584
$ws.Range("H103").Value = 27400
$ws.Range("J103").Value = 27400
$ws.Range("L103").Value = 27400
$ws.Range("N103").Value = -29744
$ws.Range("H111").Value = 17864.334
$ws.Range("J111").Value = 17864.334
$ws.Range("L111").Value = 17864.334
$ws.Range("N111").Value = -23998.334
$ws.Range("H126").Value = 1264.2858
$ws.Range("I126").Value = 885.7143
$ws.Range("J126").Value = 1642.8572
$ws.Range("K126").Value = 2657.1429
$ws.Range("L126").Value = 4928.571599999999
$ws.Range("M126").Value = -187.1428999999998
$ws.Range("N126").Value = -9868.571599999999
$ws.Range("H132").Value = 2182.4
$ws.Range("I132").Value = 1413.7778
$ws.Range("J132").Value = 3335.3333
$ws.Range("K132").Value = 4241.3334
$ws.Range("L132").Value = 10005.9999
$ws.Range("M132").Value = -1711.3334
$ws.Range("N132").Value = -15065.9999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10850.286
$ws.Range("I132").Value = 12355
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 37065
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -34535
$ws.Range("N132").Value = -21059

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3314.6562
$ws.Range("I132").Value = 4592.8
$ws.Range("K132").Value = 13778.4
$ws.Range("M132").Value = -11248.4
